$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue 'D2' '40.968.02'
Set-TextValue 'E2' '  +2.80%  '
Set-TextValue 'D3' '2.239.14'
Set-TextValue 'E3' '  +1.59%  '
Set-TextValue 'E4' '  +0.07%  '
Set-TextValue 'D5' '301.71'
Set-TextValue 'E5' '  +3.06%  '
Set-TextValue 'D6' '90.38'
Set-TextValue 'E6' '  +3.96%  '
Set-TextValue 'E7' '  +2.00%  '
Set-TextValue 'E8' '  +0.15%  '
Set-TextValue 'D9' '0.480'
Set-TextValue 'E9' '  +1.98%  '
Set-TextValue 'D10' '53.50'
Set-TextValue 'E10' '  +8.10%  '
Set-TextValue 'E11' '  +6.19%  '
Set-TextValue 'E12' '  +2.09%  '
Set-TextValue 'E13' '  +3.42%  '
Set-TextValue 'D14' '6.52'
Set-TextValue 'E14' '  +1.69%  '
Set-TextValue 'D15' '2.587.93'
Set-TextValue 'E15' '  +1.72%  '
Set-TextValue 'D16' '13.99'
Set-TextValue 'E16' '  +2.24%  '
Set-TextValue 'D17' '2.274.96'
Set-TextValue 'E17' '  +3.01%  '
Set-TextValue 'D18' '0.746'
Set-TextValue 'E18' '  +3.15%  '
Set-TextValue 'D19' '40.914.43'
Set-TextValue 'E19' '  +2.86%  '
Set-TextValue 'D20' '11.76'
Set-TextValue 'E20' '  +4.04%  '
Set-TextValue 'D21' '0.0₃0897'
Set-TextValue 'E21' '  +1.78%  '
Set-TextValue 'D22' '5.83'
Set-TextValue 'E22' '  +1.76%  '
Set-TextValue 'D23' '66.52'
Set-TextValue 'E23' '  +2.11%  '
Set-TextValue 'D24' '239.95'
Set-TextValue 'E24' '  +1.52%  '
Set-TextValue 'D25' '2.54'
Set-TextValue 'E25' '  +4.11%  '
Set-TextValue 'E26' '  +0.02%  '
Set-TextValue 'D27' '1.83'
Set-TextValue 'E27' '  +2.07%  '
Set-TextValue 'D28' '23.58'
Set-TextValue 'E28' '  +5.49%  '
Set-TextValue 'B29' 'Toncoin'
Set-TextValue 'C29' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D29' '2.17'
Set-TextValue 'E29' '  +0.54%  '
Set-TextValue 'B30' 'Cosmos'
Set-TextValue 'C30' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D30' '9.53'
Set-TextValue 'E30' '  +4.29%  '
Set-TextValue 'D31' '157.82'
Set-TextValue 'E31' '  +1.74%  '
Set-TextValue 'D32' '32.92'
Set-TextValue 'E32' '  +4.39%  '
Set-TextValue 'E33' '  +0.05%  '
Set-TextValue 'D34' '5.13'
Set-TextValue 'E34' '  +5.28%  '
Set-TextValue 'E35' '  +2.81%  '
Set-TextValue 'E36' '  +6.24%  '
Set-TextValue 'E37' '  +1.21%  '
Set-TextValue 'E38' '  +2.87%  '
Set-TextValue 'D39' '16.39'
Set-TextValue 'E39' '  +6.48%  '
Set-TextValue 'D40' '0.102'
Set-TextValue 'E40' '  +5.18%  '
Set-TextValue 'E41' '  +5.98%  '
Set-TextValue 'D42' '3.89'
Set-TextValue 'E42' '  +4.52%  '
Set-TextValue 'D43' '2.069.95'
Set-TextValue 'E43' '  -2.44%  '
Set-TextValue 'D44' '19.80'
Set-TextValue 'E44' '  +12.28%  '
Set-TextValue 'D45' '0.0274'
Set-TextValue 'E45' '  +3.27%  '
Set-TextValue 'E47' '  +10.93%  '
Set-TextValue 'E48' '  -3.71%  '
Set-TextValue 'D49' '2.460.95'
Set-TextValue 'E49' '  +1.90%  '
Set-TextValue 'E50' '  +2.38%  '
Set-TextValue 'E51' '  +3.57%  '
